$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -600
$ws.Range("B10").Value = -300
$ws.Range("C10").Value = "Shrautasutras and Grihyasutras."

$ws.Range("A11").Value = -600
$ws.Range("C11").Value = "Mahavira and the Buddha as per literature."

$ws.Range("A12").Value = -500
$ws.Range("B12").Value = -200
$ws.Range("C12").Value = "Dharmasutras."

$ws.Range("A13").Value = -450
$ws.Range("C13").Value = "Grammar of Panini."

$ws.Range("A14").Value = -500
$ws.Range("C14").Value = "Mahavira and the Buddha in the context of archaeology."

$ws.Range("A15").Value = -326
$ws.Range("C15").Value = "Alexander’s invasion."

$ws.Range("A16").Value = -322
$ws.Range("C16").Value = "Accession of Chandragupta Maurya."

$ws.Range("A17").Value = -300
$ws.Range("C17").Value = "Decipherable writing in India."

$ws.Range("A18").Value = -57
$ws.Range("C18").Value = "Vikrama Samvat."

$ws.Range("A19").Value = -100
$ws.Range("C19").Value = "Hathigumpha inscription of Kharavela of Kalinga."

$ws.Range("A20").Value = -100
$ws.Range("C20").Value = "The earliest Pali Buddhist texts compiled in Sri Lanka."

$ws.Range("A21").Value = 100
$ws.Range("C21").Value = "The Arthashastra of Kautilya finally compiled."

$ws.Range("A22").Value = 78
$ws.Range("C22").Value = "Start of Shaka Samvat."

$ws.Range("A23").Value = 80
$ws.Range("B23").Value = 115
$ws.Range("C23").Value = "The Periplus of the Erythrean Sea."

$ws.Range("A24").Value = 150
$ws.Range("C24").Value = "Ptolemy’s Geography."

$ws.Range("A25").Value = 319
$ws.Range("C25").Value = "Start of the Gupta era."

$ws.Range("A26").Value = 400
$ws.Range("C26").Value = "Mahabharata, Ramayana, and major Puranas finally compiled."

$ws.Range("A27").Value = 400
$ws.Range("C27").Value = "Earliest Indian manuscript found in Central Asia."

$ws.Range("A28").Value = 500
$ws.Range("C28").Value = "Fa-hsien comes to India."

$ws.Range("A29").Value = 600
$ws.Range("C29").Value = "The Prakrit Jaina texts finally compiled in Valabhi."

$ws.Range("A30").Value = 700
$ws.Range("C30").Value = "Hsuan Tsang’s visit. Harshacharita by Banabhatta."

$ws.Range("A31").Value = 1100
$ws.Range("C31").Value = "Mushika Vamsha by Atula."

$ws.Range("A32").Value = 1100
$ws.Range("B32").Value = 1200
$ws.Range("C32").Value = "Vikramankadevacharita by Bilhana."

$ws.Range("A33").Value = 1200
$ws.Range("C33").Value = "Ramacharita by Sandhyakara Nandi. Rajatarangini by Kalhana."

$ws.Range("A34").Value = 1837
$ws.Range("C34").Value = "Ashokan inscriptions first deciphered by James Prinsep."

$ws.PageSetup.Orientation = 1

$ws.Range("B34").Select()
